$wb = $excel.ActiveWorkbook

# Rows whose "Ready for handoff" generation/handoff timestamps move forward
# as part of regenerating the handoff report (all of them previously shared
# the same generated-at timestamp; row 12 is a different batch and is left
# untouched).
$rows = @(8, 9, 10, 11, 13, 14)

# --- "Overview" sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-17 06:18:10"
}

# --- "zh-cn" sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-17 06:17:59"
}

# --- "de-de" sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-17 06:18:10"
}
